$d = $word.ActiveDocument

$replacements = @(
    @{Old="71÷9="; New="92÷2="},
    @{Old="50÷6="; New="20÷6="},
    @{Old="81÷8="; New="75÷7="},
    @{Old="46÷2="; New="63÷6="},
    @{Old="21÷8="; New="32÷5="},
    @{Old="46÷5="; New="49÷7="},
    @{Old="95÷8="; New="69÷6="},
    @{Old="82÷4="; New="89÷4="},
    @{Old="54÷9="; New="83÷7="},
    @{Old="90÷2="; New="14÷5="},
    @{Old="25÷5="; New="25÷8="},
    @{Old="97÷8="; New="11÷8="},
    @{Old="83÷8="; New="77÷7="},
    @{Old="81÷4="; New="10÷2="},
    @{Old="49÷8="; New="43÷7="},
    @{Old="87÷8="; New="78÷5="},
    @{Old="73÷3="; New="57÷5="},
    @{Old="11÷2="; New="41÷6="},
    @{Old="57÷3="; New="65÷9="},
    @{Old="90÷7="; New="83÷4="},
    @{Old="86÷6="; New="43÷5="},
    @{Old="52÷3="; New="23÷7="},
    @{Old="68÷2="; New="61÷7="},
    @{Old="93÷5="; New="58÷9="},
    @{Old="27÷8="; New="41÷8="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
